$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Level 5 office row (row 15): base density rises from 2.5 to 24, and the
# F-column multiplier formula changes its cap from 6 to 5 and is no longer
# shared with the rows below it.
$ws.Range("E15").Value = 24
$ws.Range("F15").Formula = "=MIN(FLOOR((C15+D15)/2,1),5)-1"

# Rows 16-26 share a new F formula (cap stays at 6, parenthesised the same
# way as before but re-anchored starting at row 16).
$ws.Range("F16:F26").Formula = "=(MIN(FLOOR((C16+D16)/2,1),6)-1)"

# The "base" (E) value chain drops to 24 through row 23 (inherited
# automatically via the existing =E(prev) formulas once E15 changes),
# then steps down again to 20 starting at row 24.
$ws.Range("E24").Value = 20

# Reflect the author's last-used cell when they saved the workbook.
$ws.Range("E25").Select() | Out-Null
